# Added new login backup codes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old backup codes with the freshly generated set. The first
# three codes keep their original rows (2-4); the remaining six codes now
# occupy rows 7-12 (previously rows 13-16 held the last four old codes).
$ws.Range("A2").Value = "1CZV4VTHGTVN"
$ws.Range("A3").Value = "51K0DF5KCN34"
$ws.Range("A4").Value = "HQT8HMXSF63S"

$ws.Range("A7").Value = "5A41AVCX9PFR"
$ws.Range("A8").Value = "1YA1484DG5R7"
$ws.Range("A9").Value = "MJ6J3N01MN75"
$ws.Range("A10").Value = "B9AV6NE42R8W"
$ws.Range("A11").Value = "V2A2ZKV148W8"
$ws.Range("A12").Value = "DY9R3Z05BNS8"

# The code list shrank from 16 rows to 12 rows - clear out the now-unused
# tail so the sheet's dimension/used-range shrinks back down to A1:A12.
$ws.Range("A13:A16").ClearContents()

# Highlight the header cell with the Accent 6 theme fill.
$ws.Range("A1").Interior.ThemeColor = 10

# Move the active selection to match the resaved workbook.
$ws.Range("B6").Select() | Out-Null
